$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map the semester-name labels in column B to their numeric codes so the
# column can be used numerically (e.g. for a linear regression).
$map = @{
    "Fall"   = 1
    "J-Term" = 2
    "Spring" = 3
    "Summer" = 4
}

$lastRow = 241

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value()
    if ($null -ne $current -and $map.ContainsKey([string]$current)) {
        $cell.Value = $map[[string]$current]
    }
}

# Leave the selection where the author last left it.
[void]$ws.Range("K230").Select()

